# Daily attendance processing - 2026-01-03 07:09:13
# Applies the scraped OOXML diff to the active workbook via Excel COM-interop.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Column I ("Status") gets wider: 10 -> 14 characters.
#    The engine pads ColumnWidth by ~0.8333 on round-trip, so compensate so
#    the persisted <col width="..."/> lands on an exact 14.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 13.1666667

# ---------------------------------------------------------------------------
# 2) Every "System, dnasr281@gmail.com" note in column G (Recorded By) has
#    its two comma-separated parts swapped to "dnasr281@gmail.com, System".
#    Scan the used range instead of hard-coding rows, in case of drift.
# ---------------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count
$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}

# ---------------------------------------------------------------------------
# 3) Class-statistics block (top-right mini table): Missing/Pending counts.
# ---------------------------------------------------------------------------
$ws.Range("L7").Value = 12   # Missing Sessions: 0 -> 12
$ws.Range("L8").Value = 96   # Pending Sessions: 108 -> 96

# ---------------------------------------------------------------------------
# 4) Per-group breakdown rows 15-26: one more "missing" (P) and one fewer
#    "pending" (Q) session recorded against the running totals.
# ---------------------------------------------------------------------------
for ($r = 15; $r -le 26; $r++) {
    $ws.Cells.Item($r, 16).Value = 1   # column P
    $ws.Cells.Item($r, 17).Value = 8   # column Q
}

# ---------------------------------------------------------------------------
# 5) The twelve per-session "B1-N" summary rows flip their Status label from
#    "Pending" to "Not Recorded" (fill/format stay the same yellow style).
# ---------------------------------------------------------------------------
$summaryRows = @(18, 38, 58, 78, 97, 116, 135, 154, 173, 193, 213, 233)
foreach ($r in $summaryRows) {
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
}

Write-Host "Applied daily attendance processing edits."
